$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 82, shifting rows 82:181 down to 83:182
$ws.Rows.Item(82).Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)

# Populate the new row 82 with the new cherry-price record
$ws.Range("A82").Value = 5
$ws.Range("B82").Value = "Macroferia Regional de Talca"
$ws.Range("C82").Value = "Maule"
$ws.Range("D82").Value = 44589
$ws.Range("D82").NumberFormat = $ws.Range("D83").NumberFormat
$ws.Range("E82").Value = 7
$ws.Range("F82").Value = "Fruta"
$ws.Range("G82").Value = 100103
$ws.Range("H82").Value = "Frutos de hueso (carozo)"
$ws.Range("I82").Value = 100103001
$ws.Range("J82").Value = "Cereza"
$ws.Range("K82").Value = "Santina"
$ws.Range("L82").Value = "Primera"
$ws.Range("M82").Value = 150
$ws.Range("N82").Value = 3500
$ws.Range("O82").Value = 3500
$ws.Range("P82").Value = 3500
$ws.Range("Q82").Value = "`$/bandeja 5 kilos"
$ws.Range("R82").Value = "Provincia de Curicó"
$ws.Range("S82").Value = 700
$ws.Range("T82").Value = 5
